$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures to several
# leve rows across multiple sheets, as produced by the scheduled price-update
# runner. For each affected row we update the price/profit columns (H, I, J,
# K, L, M, N) to their new values; cells that no longer apply are cleared,
# and cells that become newly applicable are written.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 49425
$ws.Range("J134").Value = 49425
$ws.Range("L134").Value = 49425
$ws.Range("N134").Value = -59565

$ws.Range("H137").Value = 1464.5714
$ws.Range("I137").Value = 975.5
$ws.Range("J137").Value = 2116.6667
$ws.Range("K137").Value = 2926.5
$ws.Range("L137").Value = 6350.000100000001
$ws.Range("M137").Value = -376.5
$ws.Range("N137").Value = -11450.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 761.4286
$ws.Range("I74").Value = 664.8570999999999
$ws.Range("J74").Value = 954.5714
$ws.Range("K74").Value = 664.8570999999999
$ws.Range("L74").Value = 954.5714
$ws.Range("M74").Value = 209.1429000000001
$ws.Range("N74").Value = -2702.5714

$ws.Range("H77").Value = 761.4286
$ws.Range("I77").Value = 664.8570999999999
$ws.Range("J77").Value = 954.5714
$ws.Range("K77").Value = 3324.2855
$ws.Range("L77").Value = 4772.857
$ws.Range("M77").Value = 1043.7145
$ws.Range("N77").Value = -13508.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9453.429
$ws.Range("I3").Value = 11679.818
$ws.Range("K3").Value = 11679.818
$ws.Range("M3").Value = -11565.818

$ws.Range("H64").Value = 450
$ws.Range("J64").Value = 462.2857
$ws.Range("L64").Value = 462.2857
$ws.Range("N64").Value = -912.2857

$ws.Range("H67").Value = 450
$ws.Range("J67").Value = 462.2857
$ws.Range("L67").Value = 462.2857
$ws.Range("N67").Value = -2022.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 65999.664
$ws.Range("J20").Value = 65999.664
$ws.Range("L20").Value = 65999.664
$ws.Range("N20").Value = -66471.664

$ws.Range("H30").Value = 65999.664
$ws.Range("J30").Value = 65999.664
$ws.Range("L30").Value = 65999.664
$ws.Range("N30").Value = -66181.664

$ws.Range("H31").Value = 2402.2222
$ws.Range("I31").Value = 1424.3334
$ws.Range("J31").Value = 3184.5334
$ws.Range("K31").Value = 1424.3334
$ws.Range("L31").Value = 3184.5334
$ws.Range("M31").Value = -1129.3334
$ws.Range("N31").Value = -3774.5334

$ws.Range("H34").Value = 2402.2222
$ws.Range("I34").Value = 1424.3334
$ws.Range("J34").Value = 3184.5334
$ws.Range("K34").Value = 1424.3334
$ws.Range("L34").Value = 3184.5334
$ws.Range("M34").Value = -1222.3334
$ws.Range("N34").Value = -3588.5334

$ws.Range("H58").Value = 1358.5807
$ws.Range("I58").Value = 652.1579
$ws.Range("J58").Value = 2477.0833
$ws.Range("K58").Value = 652.1579
$ws.Range("L58").Value = 2477.0833
$ws.Range("M58").Value = -449.1579
$ws.Range("N58").Value = -2883.0833

$ws.Range("H128").Value = 65999.664
$ws.Range("J128").Value = 65999.664
$ws.Range("L128").Value = 65999.664
$ws.Range("N128").Value = -75959.664

$ws.Range("H132").Value = 3331.0557
$ws.Range("I132").Value = 2311
$ws.Range("J132").Value = 4351.1113
$ws.Range("K132").Value = 6933
$ws.Range("L132").Value = 13053.3339
$ws.Range("M132").Value = -4403
$ws.Range("N132").Value = -18113.3339

$ws.Range("H136").Value = 1358.5807
$ws.Range("I136").Value = 652.1579
$ws.Range("J136").Value = 2477.0833
$ws.Range("K136").Value = 1956.4737
$ws.Range("L136").Value = 7431.249899999999
$ws.Range("M136").Value = 593.5263
$ws.Range("N136").Value = -12531.2499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 357.35
$ws.Range("I23").Value = 259.5
$ws.Range("K23").Value = 778.5
$ws.Range("M23").Value = -543.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2395.6086
$ws.Range("I80").Value = 2395.6086
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2395.6086
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1397.6086
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 2395.6086
$ws.Range("I83").Value = 2395.6086
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 11978.043
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -6986.043
$ws.Range("N83").ClearContents()

$ws.Range("H92").Value = 19999.5
$ws.Range("J92").Value = 19999.5
$ws.Range("L92").Value = 19999.5
$ws.Range("N92").Value = -23743.5

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H97").Value = 2555.5557
$ws.Range("I97").Value = 2555.5557
$ws.Range("K97").Value = 2555.5557
$ws.Range("M97").Value = -2059.5557

$ws.Range("H98").Value = 23910.75
$ws.Range("J98").Value = 23910.75
$ws.Range("L98").Value = 23910.75
$ws.Range("N98").Value = -29900.75

$ws.Range("H99").Value = 4422.125
$ws.Range("I99").Value = 4422.125
$ws.Range("K99").Value = 4422.125
$ws.Range("M99").Value = -2176.125

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 12425.714
$ws.Range("J135").Value = 12425.714
$ws.Range("L135").Value = 12425.714
$ws.Range("N135").Value = -22565.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3714.4
$ws.Range("I136").Value = 2335.4119
$ws.Range("J136").Value = 5517.6924
$ws.Range("K136").Value = 7006.2357
$ws.Range("L136").Value = 16553.0772
$ws.Range("M136").Value = -4456.2357
$ws.Range("N136").Value = -21653.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 62452.824
$ws.Range("I81").Value = 201799.8
$ws.Range("J81").Value = 4391.5835
$ws.Range("K81").Value = 403599.6
$ws.Range("L81").Value = 8783.166999999999
$ws.Range("M81").Value = -402538.6
$ws.Range("N81").Value = -10905.167

$ws.Range("H84").Value = 62452.824
$ws.Range("I84").Value = 201799.8
$ws.Range("J84").Value = 4391.5835
$ws.Range("K84").Value = 2017998
$ws.Range("L84").Value = 43915.835
$ws.Range("M84").Value = -2012694
$ws.Range("N84").Value = -54523.835

$ws.Range("H132").Value = 14287235
$ws.Range("I132").Value = 20001336
$ws.Range("J132").Value = 1983.8
$ws.Range("K132").Value = 60004008
$ws.Range("L132").Value = 5951.4
$ws.Range("M132").Value = -60001478
$ws.Range("N132").Value = -11011.4

$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

Write-Output "Updated 29 leve profit rows across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."